$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Safe (non-numeric-looking) text updates - set directly
$ws.Range('D2').Value = '26.002.31'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.631.70'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -2.50%  '
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('E10').Value = '  -7.25%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '1.860.32'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '1.632.00'
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('D16').Value = '25.985.60'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').Value = '0.0₃0739'
$ws.Range('E17').Value = '  -3.53%  '
$ws.Range('E18').Value = '  -3.42%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('E22').Value = '  -3.31%  '
$ws.Range('E23').Value = '  -2.92%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('E29').Value = '  -2.83%  '
$ws.Range('E30').Value = '  -1.85%  '
$ws.Range('E31').Value = '  -3.66%  '
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('E33').Value = '  -5.61%  '
$ws.Range('E34').Value = '  -2.26%  '
$ws.Range('E35').Value = '  -4.03%  '
$ws.Range('D36').Value = '1.128.17'
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  -5.43%  '
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('E39').Value = '  -4.95%  '
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('E42').Value = '  -3.19%  '
$ws.Range('D43').Value = '1.769.79'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('E44').Value = '  -5.33%  '
$ws.Range('D45').Value = '0.0₆0115'
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  -3.71%  '

# Numeric-looking text values: force Text format to preserve exact string representation,
# then restore default "Normal" style so no stray formatting is introduced.
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '214.38'
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '18.25'
$r.Style = 'Normal'
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '0.0791'
$r.Style = 'Normal'
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '0.524'
$r.Style = 'Normal'
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '190.15'
$r.Style = 'Normal'
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '9.62'
$r.Style = 'Normal'
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '6.06'
$r.Style = 'Normal'
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '1.79'
$r.Style = 'Normal'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '144.17'
$r.Style = 'Normal'
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '6.75'
$r.Style = 'Normal'
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '15.14'
$r.Style = 'Normal'
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '0.0479'
$r.Style = 'Normal'
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '3.11'
$r.Style = 'Normal'
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.858'
$r.Style = 'Normal'
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.516'
$r.Style = 'Normal'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '98.30'
$r.Style = 'Normal'
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '0.773'
$r.Style = 'Normal'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '54.76'
$r.Style = 'Normal'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '0.0526'
$r.Style = 'Normal'
